$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 823
$ws1.Range("F6").Value = 390
$ws1.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202410/SodCscpF1729479237675.jpeg"
$ws1.Range("F12").Value = 13305
$ws1.Range("F14").Value = 11
$ws1.Range("F16").Value = 5474
$ws1.Range("F17").Value = 5568
$ws1.Range("F18").Value = 39

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 823
$ws4.Range("F22").Value = 390
$ws4.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202410/SodCscpF1729479237675.jpeg"
$ws4.Range("F34").Value = 13305
$ws4.Range("F36").Value = 11
$ws4.Range("F39").Value = 5474
$ws4.Range("F40").Value = 5568
$ws4.Range("F41").Value = 39
